$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-21 Wednesday" "2025-05-22 Thursday"

Replace-Text "614÷5=" "523÷8="
Replace-Text "127÷2=" "571÷3="
Replace-Text "681÷4=" "262÷5="
Replace-Text "130÷6=" "845÷6="
Replace-Text "773÷3=" "362÷8="
Replace-Text "699÷5=" "733÷2="
Replace-Text "669÷4=" "991÷9="
Replace-Text "818÷7=" "988÷9="
Replace-Text "899÷8=" "746÷4="
Replace-Text "288÷9=" "431÷2="
Replace-Text "951÷9=" "200÷9="
Replace-Text "774÷3=" "172÷5="
Replace-Text "461÷9=" "762÷5="
Replace-Text "478÷4=" "521÷4="
Replace-Text "574÷8=" "528÷8="
Replace-Text "935÷9=" "996÷3="
Replace-Text "181÷6=" "395÷3="
Replace-Text "849÷9=" "189÷6="
Replace-Text "169÷8=" "128÷9="
Replace-Text "523÷6=" "539÷2="
Replace-Text "986÷3=" "270÷4="
Replace-Text "707÷8=" "109÷8="
Replace-Text "761÷5=" "475÷9="
Replace-Text "238÷3=" "287÷5="
Replace-Text "224÷2=" "669÷5="
